# Auto-generated: update market price / profit columns (H-N) for specific leve rows
# across multiple crafting-class sheets, per scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1562.7333
$ws.Range("I43").Value = 913.9
$ws.Range("J43").Value = 2860.4
$ws.Range("K43").Value = 913.9
$ws.Range("L43").Value = 2860.4
$ws.Range("M43").Value = -844.9
$ws.Range("N43").Value = -2998.4

$ws.Range("H55").Value = 700.9231
$ws.Range("J55").Value = 751.75
$ws.Range("L55").Value = 751.75
$ws.Range("N55").Value = -1179.75

$ws.Range("H98").Value = 1882.9584
$ws.Range("I98").Value = 853.73914
$ws.Range("J98").Value = 25555
$ws.Range("K98").Value = 853.73914
$ws.Range("L98").Value = 25555
$ws.Range("M98").Value = 644.26086
$ws.Range("N98").Value = -28551

$ws.Range("H122").Value = 1882.9584
$ws.Range("I122").Value = 853.73914
$ws.Range("J122").Value = 25555
$ws.Range("K122").Value = 2561.21742
$ws.Range("L122").Value = 76665
$ws.Range("M122").Value = -111.2174199999999
$ws.Range("N122").Value = -81565

$ws.Range("H132").Value = 1955.5161
$ws.Range("I132").Value = 1309.7727
$ws.Range("K132").Value = 3929.3181
$ws.Range("M132").Value = -1399.3181

$ws.Range("H135").Value = 4491.56
$ws.Range("I135").Value = 5950.0557
$ws.Range("J135").Value = 741.1429
$ws.Range("K135").Value = 53550.5013
$ws.Range("L135").Value = 6670.2861
$ws.Range("M135").Value = -51015.5013
$ws.Range("N135").Value = -11740.2861

$ws.Range("H141").Value = 4300.4688
$ws.Range("I141").Value = 2402.8948
$ws.Range("J141").Value = 7073.846
$ws.Range("K141").Value = 7208.6844
$ws.Range("L141").Value = 21221.538
$ws.Range("M141").Value = -2028.6844
$ws.Range("N141").Value = -31581.538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1273.8846
$ws.Range("I61").Value = 509.8125
$ws.Range("J61").Value = 2496.4
$ws.Range("K61").Value = 509.8125
$ws.Range("L61").Value = 2496.4
$ws.Range("M61").Value = -297.8125
$ws.Range("N61").Value = -2920.4

$ws.Range("H97").Value = 794.45
$ws.Range("J97").Value = 1005.5
$ws.Range("L97").Value = 1005.5
$ws.Range("N97").Value = -1997.5

$ws.Range("H132").Value = 1960.4043
$ws.Range("I132").Value = 998.2593
$ws.Range("J132").Value = 3259.3
$ws.Range("K132").Value = 2994.7779
$ws.Range("L132").Value = 9777.900000000001
$ws.Range("M132").Value = -464.7779
$ws.Range("N132").Value = -14837.9

$ws.Range("H136").Value = 1273.8846
$ws.Range("I136").Value = 509.8125
$ws.Range("J136").Value = 2496.4
$ws.Range("K136").Value = 1529.4375
$ws.Range("L136").Value = 7489.200000000001
$ws.Range("M136").Value = 1020.5625
$ws.Range("N136").Value = -12589.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1228.2106
$ws.Range("I107").Value = 1242.5834
$ws.Range("K107").Value = 1242.5834
$ws.Range("M107").Value = 677.4166

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2657.96
$ws.Range("I16").Value = 1359.9333
$ws.Range("J16").Value = 4605
$ws.Range("K16").Value = 1359.9333
$ws.Range("L16").Value = 4605
$ws.Range("M16").Value = -1072.9333
$ws.Range("N16").Value = -5179

$ws.Range("H31").Value = 1495.904
$ws.Range("I31").Value = 1103.575
$ws.Range("J31").Value = 1971.4546
$ws.Range("K31").Value = 1103.575
$ws.Range("L31").Value = 1971.4546
$ws.Range("M31").Value = -808.575
$ws.Range("N31").Value = -2561.4546

$ws.Range("H34").Value = 1495.904
$ws.Range("I34").Value = 1103.575
$ws.Range("J34").Value = 1971.4546
$ws.Range("K34").Value = 1103.575
$ws.Range("L34").Value = 1971.4546
$ws.Range("M34").Value = -901.575
$ws.Range("N34").Value = -2375.4546

$ws.Range("H99").Value = 954.6667
$ws.Range("I99").Value = 938.4
$ws.Range("J99").Value = 966.2857
$ws.Range("K99").Value = 938.4
$ws.Range("L99").Value = 966.2857
$ws.Range("M99").Value = 559.6
$ws.Range("N99").Value = -3962.2857

$ws.Range("H107").Value = 1042.1111
$ws.Range("I107").Value = 385.35715
$ws.Range("K107").Value = 385.35715
$ws.Range("M107").Value = 1534.64285

$ws.Range("H113").Value = 2657.96
$ws.Range("I113").Value = 1359.9333
$ws.Range("J113").Value = 4605
$ws.Range("K113").Value = 1359.9333
$ws.Range("L113").Value = 4605
$ws.Range("M113").Value = 810.0667000000001
$ws.Range("N113").Value = -8945

$ws.Range("H122").Value = 1053.5
$ws.Range("I122").Value = 1250
$ws.Range("J122").Value = 857
$ws.Range("K122").Value = 3750
$ws.Range("L122").Value = 2571
$ws.Range("M122").Value = -1300
$ws.Range("N122").Value = -7471

$ws.Range("H126").Value = 954.6667
$ws.Range("I126").Value = 938.4
$ws.Range("J126").Value = 966.2857
$ws.Range("K126").Value = 2815.2
$ws.Range("L126").Value = 2898.8571
$ws.Range("M126").Value = -345.1999999999998
$ws.Range("N126").Value = -7838.8571

$ws.Range("H132").Value = 2646.4167
$ws.Range("I132").Value = 2282.3
$ws.Range("J132").Value = 4467
$ws.Range("K132").Value = 6846.900000000001
$ws.Range("L132").Value = 13401
$ws.Range("M132").Value = -4316.900000000001
$ws.Range("N132").Value = -18461

$ws.Range("H134").Value = 1256.875
$ws.Range("I134").Value = 1180.129
$ws.Range("J134").Value = 1396.8235
$ws.Range("K134").Value = 3540.387
$ws.Range("L134").Value = 4190.470499999999
$ws.Range("M134").Value = -1005.387
$ws.Range("N134").Value = -9260.4705

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 813.5
$ws.Range("I132").Value = 852
$ws.Range("J132").Value = 775
$ws.Range("K132").Value = 7668
$ws.Range("L132").Value = 6975
$ws.Range("M132").Value = -5138
$ws.Range("N132").Value = -12035

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5212.636
$ws.Range("I132").Value = 6352.136
$ws.Range("J132").Value = 2933.6365
$ws.Range("K132").Value = 19056.408
$ws.Range("L132").Value = 8800.9095
$ws.Range("M132").Value = -16526.408
$ws.Range("N132").Value = -13860.9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2020.4
$ws.Range("I7").Value = 1913
$ws.Range("J7").Value = 2450
$ws.Range("K7").Value = 1913
$ws.Range("L7").Value = 2450
$ws.Range("M7").Value = -1801
$ws.Range("N7").Value = -2674

$ws.Range("H61").Value = 1281.4
$ws.Range("I61").Value = 1281.4
$ws.Range("K61").Value = 1281.4
$ws.Range("M61").Value = -1079.4

$ws.Range("H113").Value = 1281.4
$ws.Range("I113").Value = 1281.4
$ws.Range("K113").Value = 1281.4
$ws.Range("M113").Value = 888.5999999999999

$ws.Range("H126").Value = 2020.4
$ws.Range("I126").Value = 1913
$ws.Range("J126").Value = 2450
$ws.Range("K126").Value = 5739
$ws.Range("L126").Value = 7350
$ws.Range("M126").Value = -3269
$ws.Range("N126").Value = -12290

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 6061433
$ws.Range("I100").Value = 7576624.5
$ws.Range("J100").Value = 666.6667
$ws.Range("K100").Value = 15153249
$ws.Range("L100").Value = 1333.3334
$ws.Range("M100").Value = -15152708
$ws.Range("N100").Value = -2415.3334
